$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Pieza de fabricante F8GPYPieza de Dell 210-BFRQ'
$ws.Cells.Item(2, 2).Value = 'Monitor curvo Alienware 34 QD-OLED para juegos - AW3423DWF'
$ws.Cells.Item(2, 3).Value = 914.76
$ws.Cells.Item(2, 4).Value = 756
$ws.Cells.Item(2, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/monitors/alienware/aw3423dwf/media-gallery/monitor-alienware-aw3423dwf-black-gallery-10.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(3, 1).Value = 'Pieza de fabricante PM80JPieza de Dell 210-BHTK'
$ws.Cells.Item(3, 2).Value = 'Monitor gaming Dell 27 - G2724D'
$ws.Cells.Item(3, 3).Value = 217.99
$ws.Cells.Item(3, 4).Value = 180.16
$ws.Cells.Item(3, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/monitors/g-series/g2724d/media-gallery/monitor-g2724d-black-gallery-1.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(4, 1).Value = 'Pieza de fabricante 71YCFPieza de Dell 210-AZZD'
$ws.Cells.Item(4, 2).Value = 'Monitor curvo para juegos Dell de 68,58 cm (27") (S2722DGM)'
$ws.Cells.Item(4, 3).Value = 168.43
$ws.Cells.Item(4, 4).Value = 139.2
$ws.Cells.Item(4, 5).Value = '//i.dell.com/is/image/DellContent//content/dam/ss2/product-images/dell-client-products/peripherals/monitors/s-series/s2722dgm/media-gallery/s2722dgm_cfp_00000ff090_bk.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(5, 1).Value = 'Pieza de fabricante HF0CGPieza de Dell 210-AZZE'
$ws.Cells.Item(5, 2).Value = 'Monitor curvo para juegos Dell 34 – S3422DWG'
$ws.Cells.Item(5, 3).Value = 364.98
$ws.Cells.Item(5, 4).Value = 301.64
$ws.Cells.Item(5, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/monitors/s-series/s3422dwg/media-gallery/s3422dwg_xfp_01_bk.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(6, 1).Value = 'Pieza de fabricante KM03TPieza de Dell 210-BLHH'
$ws.Cells.Item(6, 2).Value = 'Monitor gaming QD-OLED Alienware de 68,58 cm (27) a 360 Hz - AW2725DF'
$ws.Cells.Item(6, 3).Value = 698.99
$ws.Cells.Item(6, 4).Value = 577.68
$ws.Cells.Item(6, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/monitors/alienware/aw2725df/media-gallery/monitor-alienware-2725df-black-gallery-1-mg.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(7, 1).Value = 'Pieza de fabricante 7KYY6Pieza de Dell 545-BBFT'
$ws.Cells.Item(7, 2).Value = 'Auriculares gaming inalámbricos Alienware Pro'
$ws.Cells.Item(7, 3).Value = 255
$ws.Cells.Item(7, 4).Value = 210.74
$ws.Cells.Item(7, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/headphones/aw-pro-wireless-headset/media-gallery/lunar-light/headset-aw-pro-wh-gallery-1.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(8, 1).Value = 'Pieza de fabricante R6FWGPieza de Dell 460-BDIC'
$ws.Cells.Item(8, 2).Value = 'Mochila de trabajo Alienware Horizon'
$ws.Cells.Item(8, 3).Value = 59.99
$ws.Cells.Item(8, 4).Value = 49.58
$ws.Cells.Item(8, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/peripherals/alienware/backpack/aw523p/aw523p-mg-1-5000.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(9, 1).Value = 'Pieza de fabricante 83XJYPieza de Dell 545-BBFS'
$ws.Cells.Item(9, 2).Value = 'Auriculares gaming inalámbricos Alienware Pro'
$ws.Cells.Item(9, 3).Value = 255
$ws.Cells.Item(9, 4).Value = 210.74
$ws.Cells.Item(9, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/headphones/aw-pro-wireless-headset/media-gallery/dark-side-of-the-moon/headset-aw-pro-bk-gallery-1.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(10, 1).Value = 'Pieza de fabricante 8VY4KPieza de Dell 210-BDXS'
$ws.Cells.Item(10, 2).Value = 'Monitor Dell 32 UHD 4K para juegos - G3223Q'
$ws.Cells.Item(10, 3).Value = 519.96
$ws.Cells.Item(10, 4).Value = 429.72
$ws.Cells.Item(10, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/monitors/g-series/g3223q/media-gallery/monitor-g3223q-black-gallery-1.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(11, 1).Value = 'Pieza de fabricante J1PTNPieza de Dell 210-AZZH'
$ws.Cells.Item(11, 2).Value = 'Monitor curvo para juegos Dell de 81,28 cm (32") (S3222DGM)'
$ws.Cells.Item(11, 3).Value = 344.85
$ws.Cells.Item(11, 4).Value = 285
$ws.Cells.Item(11, 5).Value = '//i.dell.com/is/image/DellContent//content/dam/ss2/product-images/dell-client-products/peripherals/monitors/s-series/s3222dgm/media-gallery/s3222dgm_cfp_00000ff090_bk.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(12, 1).Value = 'Pieza de fabricante MPJG7Pieza de Dell 210-BLLV'
$ws.Cells.Item(12, 2).Value = 'Monitor gaming QD-OLED 4K Alienware de 81,28 cm (32") - AW3225QF'
$ws.Cells.Item(12, 3).Value = 849.96
$ws.Cells.Item(12, 4).Value = 702.45
$ws.Cells.Item(12, 5).Value = '//i.dell.com/is/image/DellContent/content/dam/ss2/product-images/dell-client-products/peripherals/monitors/aw-series/aw3225qf/media-gallery/monitor-alienware-aw3225qf-white-gallery-1.psd?qlt=90,0&op_usm=1.75,0.3,2,0&resMode=sharp&pscan=auto&fmt=png-alpha&hei=500'

$ws.Cells.Item(13, 1).Value = 'Pieza de fabricante CNH4JPieza de Dell 460-BCYY'
$ws.Cells.Item(13, 2).Value = 'Mochila Dell Gaming 17'
$ws.Cells.Item(13, 3).Value = 58.36
$ws.Cells.Item(13, 4).Value = 48.23
$ws.Cells.Item(13, 5).Value = '//snpi.dell.com/snp/images/products/large/es-es~460-BCYY_V1%20-%20Copy/460-BCYY_V1%20-%20Copy.jpg'
